$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.848.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.72%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.524.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.95%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'606.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.92%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'197.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +6.03%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.47%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.09%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -7.50%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -0.71%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'53.80"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.33%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -2.47%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'9.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.39%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'4.083.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.05%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'598.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.91%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'70.036.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.97%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'19.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.43%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'12.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.27%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.528.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.15%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +1.37%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.993"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.54%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'18.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +6.66%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +4.50%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'102.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.62%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'4.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.08%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +5.49%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'10.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.61%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'9.63"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.94%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.66%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +10.62%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'7.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.72%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +0.82%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -0.26%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'63.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.00%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.0₃0858"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +10.38%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'3.734.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +3.82%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'3.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -3.24%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.15%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +0.05%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.393"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.79%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'36.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.20%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'488.87"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -5.63%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.134"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.39%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -1.70%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -3.28%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -1.96%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -1.55%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +0.38%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'8.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -3.26%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.000249"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.66%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +11.12%  "
$ws.Range("E51").Style = "Normal"
